# Refresh Halicarnassus market-board profit data across all Sheets.
# Each worksheet's H/I/J/K/L/M/N (price + profit) columns are updated
# to the latest scheduled-runner snapshot. A handful of rows lose their
# NQ or HQ profit cell entirely when that side's price becomes 0 in the
# refreshed snapshot (mirrors the source generator, which omits the
# profit cell rather than writing a 0/undefined profit).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 25
$ws.Range("I11").Value = 25
$ws.Range("K11").Value = 25
$ws.Range("M11").Value = 115
$ws.Range("H12").Value = 100
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 100
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -440
$ws.Range("H15").Value = 4151.3335
$ws.Range("I15").Value = 4151.3335
$ws.Range("K15").Value = 12454.0005
$ws.Range("M15").Value = -12285.0005
$ws.Range("H38").Value = 1514.6666
$ws.Range("J38").Value = 2376.1428
$ws.Range("L38").Value = 7128.428400000001
$ws.Range("N38").Value = -7872.428400000001
$ws.Range("H43").Value = 795.5
$ws.Range("I43").Value = 590
$ws.Range("K43").Value = 590
$ws.Range("M43").Value = -521
$ws.Range("H58").Value = 2330.1333
$ws.Range("J58").Value = 3154.2727
$ws.Range("L58").Value = 9462.8181
$ws.Range("N58").Value = -9762.8181
$ws.Range("H92").Value = 240
$ws.Range("I92").Value = 240
$ws.Range("K92").Value = 240
$ws.Range("M92").Value = 1008
$ws.Range("H94").Value = 3068.6155
$ws.Range("I94").Value = 1899.2727
$ws.Range("J94").Value = 9500
$ws.Range("K94").Value = 1899.2727
$ws.Range("L94").Value = 9500
$ws.Range("M94").Value = -1448.2727
$ws.Range("N94").Value = -10402
$ws.Range("H107").Value = 227.85715
$ws.Range("I107").Value = 213.45454
$ws.Range("K107").Value = 213.45454
$ws.Range("M107").Value = 1706.54546
$ws.Range("H141").Value = 8333
$ws.Range("I141").Value = 7499.5
$ws.Range("K141").Value = 22498.5
$ws.Range("M141").Value = -17318.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20001708
$ws.Range("I32").Value = 2135.5
$ws.Range("K32").Value = 2135.5
$ws.Range("M32").Value = -1848.5
$ws.Range("H45").Value = 4013.2
$ws.Range("I45").Value = 2611
$ws.Range("K45").Value = 2611
$ws.Range("M45").Value = -2234
$ws.Range("H74").Value = 2654.9333
$ws.Range("I74").Value = 2253.182
$ws.Range("K74").Value = 2253.182
$ws.Range("M74").Value = -1379.182
$ws.Range("H77").Value = 2654.9333
$ws.Range("I77").Value = 2253.182
$ws.Range("K77").Value = 11265.91
$ws.Range("M77").Value = -6897.91
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1651.1666
$ws.Range("I20").Value = 2001.75
$ws.Range("J20").Value = 950
$ws.Range("K20").Value = 2001.75
$ws.Range("L20").Value = 950
$ws.Range("M20").Value = -1754.75
$ws.Range("N20").Value = -1444
$ws.Range("H105").Value = 22185616
$ws.Range("I105").Value = 22185616
$ws.Range("K105").Value = 22185616
$ws.Range("M105").Value = -22183869

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 58882.5
$ws.Range("J51").Value = 58882.5
$ws.Range("L51").Value = 58882.5
$ws.Range("N51").Value = -60354.5
$ws.Range("H58").Value = 3265.4443
$ws.Range("I58").Value = 1628.4286
$ws.Range("J58").Value = 8995
$ws.Range("K58").Value = 1628.4286
$ws.Range("L58").Value = 8995
$ws.Range("M58").Value = -1425.4286
$ws.Range("N58").Value = -9401
$ws.Range("H59").Value = 42441.25
$ws.Range("J59").Value = 42441.25
$ws.Range("L59").Value = 42441.25
$ws.Range("N59").Value = -44731.25
$ws.Range("H60").Value = 22716.555
$ws.Range("I60").Value = 3669.1428
$ws.Range("J60").Value = 89382.5
$ws.Range("K60").Value = 3669.1428
$ws.Range("L60").Value = 89382.5
$ws.Range("M60").Value = -3158.1428
$ws.Range("N60").Value = -90404.5
$ws.Range("H61").Value = 58882.5
$ws.Range("J61").Value = 58882.5
$ws.Range("L61").Value = 58882.5
$ws.Range("N61").Value = -59578.5
$ws.Range("H86").Value = 4252.5
$ws.Range("I86").Value = 4252.5
$ws.Range("K86").Value = 4252.5
$ws.Range("M86").Value = -3129.5
$ws.Range("H89").Value = 4252.5
$ws.Range("I89").Value = 4252.5
$ws.Range("K89").Value = 21262.5
$ws.Range("M89").Value = -15646.5
$ws.Range("H122").Value = 408.66666
$ws.Range("I122").Value = 370.4
$ws.Range("K122").Value = 1111.2
$ws.Range("M122").Value = 1338.8
$ws.Range("H132").Value = 2410.4443
$ws.Range("I132").Value = 2410.4443
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7231.3329
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4701.3329
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 3258.08
$ws.Range("I134").Value = 2548.6667
$ws.Range("J134").Value = 6982.5
$ws.Range("K134").Value = 7646.000100000001
$ws.Range("L134").Value = 20947.5
$ws.Range("M134").Value = -5111.000100000001
$ws.Range("N134").Value = -26017.5
$ws.Range("H136").Value = 3265.4443
$ws.Range("I136").Value = 1628.4286
$ws.Range("J136").Value = 8995
$ws.Range("K136").Value = 4885.2858
$ws.Range("L136").Value = 26985
$ws.Range("M136").Value = -2335.2858
$ws.Range("N136").Value = -32085

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 8439.9
$ws.Range("I138").Value = 3500
$ws.Range("J138").Value = 9674.875
$ws.Range("K138").Value = 10500
$ws.Range("L138").Value = 29024.625
$ws.Range("M138").Value = -5360
$ws.Range("N138").Value = -39304.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H113").Value = 6994.9165
$ws.Range("I113").Value = 985
$ws.Range("K113").Value = 985
$ws.Range("M113").Value = 1185
$ws.Range("H132").Value = 145664.14
$ws.Range("I132").Value = 169400.17
$ws.Range("K132").Value = 508200.51
$ws.Range("M132").Value = -505670.51

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1502.0588
$ws.Range("I93").Value = 1402.7333
$ws.Range("K93").Value = 1402.7333
$ws.Range("M93").Value = -154.7333000000001
$ws.Range("H100").Value = 4713
$ws.Range("I100").Value = 1307.1666
$ws.Range("K100").Value = 1307.1666
$ws.Range("M100").Value = -766.1666
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H122").Value = 2999.3845
$ws.Range("I122").Value = 2999.3635
$ws.Range("K122").Value = 8998.0905
$ws.Range("M122").Value = -6548.0905
$ws.Range("H130").Value = 34999.5
$ws.Range("J130").Value = 34999.5
$ws.Range("L130").Value = 34999.5
$ws.Range("N130").Value = -45039.5
$ws.Range("H140").Value = 125000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 22333.334
$ws.Range("J112").Value = 22333.334
$ws.Range("L112").Value = 22333.334
$ws.Range("N112").Value = -25287.334
$ws.Range("H122").Value = 3918.9
$ws.Range("I122").Value = 2537.8
$ws.Range("K122").Value = 7613.400000000001
$ws.Range("M122").Value = -5163.400000000001
$ws.Range("H125").Value = 42500
$ws.Range("J125").Value = 42500
$ws.Range("L125").Value = 42500
$ws.Range("N125").Value = -52340
$ws.Range("H126").Value = 4005.6
$ws.Range("J126").Value = 5938.8
$ws.Range("L126").Value = 17816.4
$ws.Range("N126").Value = -22756.4
$ws.Range("H132").Value = 1030.4445
$ws.Range("I132").Value = 1059
$ws.Range("J132").Value = 545
$ws.Range("K132").Value = 3177
$ws.Range("L132").Value = 1635
$ws.Range("M132").Value = -647
$ws.Range("N132").Value = -6695
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 2142.2222
$ws.Range("J136").Value = 3497.9167
$ws.Range("L136").Value = 10493.7501
$ws.Range("N136").Value = -15593.7501
